$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.676915526390076
$ws.Range("B1").Value = 2.823069095611572
$ws.Range("C1").Value = 5.8324294090271
$ws.Range("D1").Value = 2.308043956756592
$ws.Range("E1").Value = 0.8184143900871277
